$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("D5").Value = "2016-26-11 14:26:39"
$wsZhCn.Range("E5").Value = "2016-03-11 14:26:36"
$wsDeDe.Range("E5").Value = "2016-03-11 14:26:39"
